# Feature #241 verification run: the TEST-FEATURE-139-SORTIE placeholder
# sortie row is removed from the report, all subsequent sortie rows shift
# up by one, the report header (generated timestamp + total sortie count)
# is refreshed, and the CUI footer row follows the data up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete test row (row 8: TEST-FEATURE-139-SORTIE). This
# shifts every row below it (CRIIS-SORTIE-004, -002, -001, -003 and the
# CUI footer) up by one row automatically.
$ws.Rows(8).Delete()

# Refresh the report header text to reflect the regenerated report.
$ws.Range("A4").Value = "Generated: 2026-01-20 09:56:08Z"
$ws.Range("A5").Value = "Total Sorties: 4"
